$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainText($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $value
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

function Set-PlainValue($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

Set-PlainValue 'D2' '42.279.41'
Set-PlainValue 'E2' '  -2.57%  '
Set-PlainValue 'D3' '2.222.23'
Set-PlainValue 'E3' '  -1.98%  '
Set-PlainText 'D4' '1.00'
Set-PlainValue 'E4' '  +0.05%  '
Set-PlainText 'D5' '109.50'
Set-PlainValue 'E5' '  -8.12%  '
Set-PlainText 'D6' '296.45'
Set-PlainValue 'E6' '  +11.74%  '
Set-PlainText 'D7' '0.620'
Set-PlainValue 'E7' '  -4.02%  '
Set-PlainValue 'E8' '  -0.28%  '
Set-PlainText 'D9' '0.600'
Set-PlainValue 'E9' '  -2.99%  '
Set-PlainText 'D10' '43.87'
Set-PlainValue 'E10' '  -7.44%  '
Set-PlainValue 'E11' '  -3.16%  '
Set-PlainText 'D12' '54.21'
Set-PlainValue 'E12' '  -0.07%  '
Set-PlainText 'D13' '8.80'
Set-PlainValue 'E13' '  -4.06%  '
Set-PlainText 'D14' '0.999'
Set-PlainValue 'E14' '  +10.85%  '
Set-PlainValue 'E15' '  -2.64%  '
Set-PlainText 'D16' '15.09'
Set-PlainValue 'E16' '  -1.98%  '
Set-PlainValue 'D17' '2.551.98'
Set-PlainValue 'E17' '  -2.23%  '
Set-PlainValue 'D18' '2.224.09'
Set-PlainValue 'E18' '  -1.90%  '
Set-PlainValue 'D19' '42.319.56'
Set-PlainValue 'E19' '  -2.83%  '
Set-PlainText 'D20' '7.36'
Set-PlainValue 'E20' '  +7.50%  '
Set-PlainValue 'E21' '  -3.97%  '
Set-PlainText 'D22' '72.25'
Set-PlainValue 'E22' '  +0.27%  '
Set-PlainText 'D23' '3.47'
Set-PlainValue 'E23' '  +21.09%  '
Set-PlainValue 'E24' '  -3.14%  '
Set-PlainText 'D25' '227.56'
Set-PlainValue 'E25' '  -3.33%  '
Set-PlainText 'D26' '9.11'
Set-PlainValue 'E26' '  -4.17%  '
Set-PlainValue 'B27' 'Cosmos'
Set-PlainValue 'C27' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-PlainText 'D27' '11.66'
Set-PlainValue 'E27' '  -2.88%  '
Set-PlainValue 'B28' 'Dai'
Set-PlainValue 'C28' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-PlainText 'D28' '0.999'
Set-PlainValue 'E28' '  -1.76%  '
Set-PlainValue 'E29' '  -1.21%  '
Set-PlainText 'D30' '38.20'
Set-PlainValue 'E30' '  -8.74%  '
Set-PlainText 'D31' '174.01'
Set-PlainValue 'E31' '  +1.30%  '
Set-PlainText 'D32' '3.20'
Set-PlainValue 'E32' '  -5.49%  '
Set-PlainText 'D33' '20.99'
Set-PlainValue 'E33' '  -2.73%  '
Set-PlainValue 'E34' '  -1.72%  '
Set-PlainText 'D35' '5.61'
Set-PlainValue 'E35' '  -1.74%  '
Set-PlainText 'D36' '5.06'
Set-PlainValue 'E36' '  +11.19%  '
Set-PlainValue 'E37' '  +3.17%  '
Set-PlainValue 'B38' 'Stellar'
Set-PlainValue 'C38' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-PlainText 'D38' '0.126'
Set-PlainValue 'E38' '  -3.41%  '
Set-PlainValue 'B39' 'VeChain'
Set-PlainValue 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-PlainText 'D39' '0.0380'
Set-PlainValue 'E39' '  -0.39%  '
Set-PlainValue 'E40' '  -2.69%  '
Set-PlainText 'D41' '2.42'
Set-PlainValue 'E41' '  -5.05%  '
Set-PlainText 'D42' '71.91'
Set-PlainValue 'E42' '  -2.79%  '
Set-PlainText 'D43' '0.233'
Set-PlainValue 'E43' '  -1.66%  '
Set-PlainValue 'E44' '  +0.07%  '
Set-PlainText 'D45' '12.56'
Set-PlainValue 'E45' '  -9.41%  '
Set-PlainValue 'E46' '  -4.33%  '
Set-PlainText 'D47' '5.42'
Set-PlainValue 'E47' '  -6.59%  '
Set-PlainText 'D49' '103.15'
Set-PlainValue 'E49' '  +1.67%  '
Set-PlainValue 'B50' 'FraxShare'
Set-PlainValue 'C50' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-PlainText 'D50' '8.43'
Set-PlainValue 'E50' '  -1.46%  '
Set-PlainValue 'B51' 'Stacks'
Set-PlainValue 'C51' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-PlainText 'D51' '1.64'
Set-PlainValue 'E51' '  +5.98%  '
